$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.091.39"
$ws.Range("E2").Value = "'  -0.78%  "
$ws.Range("D3").Value = "'3.519.68"
$ws.Range("E3").Value = "'  +0.22%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'604.60"
$ws.Range("E5").Value = "'  -0.82%  "
$ws.Range("D6").Value = "'148.25"
$ws.Range("E6").Value = "'  -2.78%  "
$ws.Range("D7").Value = "'3.519.90"
$ws.Range("E7").Value = "'  +0.26%  "
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = "'  -1.41%  "
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "'  -0.76%  "
$ws.Range("D11").Value = "'7.83"
$ws.Range("E11").Value = "'  +3.22%  "
$ws.Range("D12").Value = "'0.424"
$ws.Range("E12").Value = "'  -1.92%  "
$ws.Range("D13").Value = "'0.0000216"
$ws.Range("E13").Value = "'  -0.66%  "
$ws.Range("D14").Value = "'4.104.94"
$ws.Range("E14").Value = "'  +0.04%  "
$ws.Range("D15").Value = "'31.67"
$ws.Range("E15").Value = "'  -2.87%  "
$ws.Range("D16").Value = "'3.511.74"
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("D17").Value = "'67.112.85"
$ws.Range("E17").Value = "'  -0.43%  "
$ws.Range("E18").Value = "'  -0.70%  "
$ws.Range("D19").Value = "'10.74"
$ws.Range("E19").Value = "'  +8.82%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("E20").Value = "'  -2.31%  "
$ws.Range("D21").Value = "'15.41"
$ws.Range("E21").Value = "'  -1.27%  "
$ws.Range("D22").Value = "'436.80"
$ws.Range("E22").Value = "'  -2.31%  "
$ws.Range("D23").Value = "'0.612"
$ws.Range("E23").Value = "'  -3.23%  "
$ws.Range("D24").Value = "'80.11"
$ws.Range("E24").Value = "'  +2.47%  "
$ws.Range("D25").Value = "'3.645.39"
$ws.Range("E25").Value = "'  -0.13%  "
$ws.Range("E26").Value = "'  +0.43%  "
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "'  -4.70%  "
$ws.Range("D28").Value = "'9.87"
$ws.Range("E28").Value = "'  -2.38%  "
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "'  -5.32%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "'  -0.23%  "
$ws.Range("D31").Value = "'1.61"
$ws.Range("E31").Value = "'  -3.17%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("D33").Value = "'0.167"
$ws.Range("E33").Value = "'  -2.06%  "
$ws.Range("D34").Value = "'25.39"
$ws.Range("E34").Value = "'  -1.35%  "
$ws.Range("D35").Value = "'3.501.97"
$ws.Range("E35").Value = "'  -0.10%  "
$ws.Range("D36").Value = "'1.80"
$ws.Range("E36").Value = "'  -3.65%  "
$ws.Range("D37").Value = "'5.92"
$ws.Range("E37").Value = "'  -4.42%  "
$ws.Range("D38").Value = "'8.04"
$ws.Range("E38").Value = "'  -0.11%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "'  -0.13%  "
$ws.Range("D41").Value = "'0.0895"
$ws.Range("E41").Value = "'  -0.11%  "
$ws.Range("D42").Value = "'170.04"
$ws.Range("E42").Value = "'  -1.79%  "
$ws.Range("B43").Value = "'Stacks"
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.10"
$ws.Range("E43").Value = "'  -9.62%  "
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'5.44"
$ws.Range("E44").Value = "'  -0.63%  "
$ws.Range("D45").Value = "'0.897"
$ws.Range("E45").Value = "'  +1.40%  "
$ws.Range("D46").Value = "'29.17"
$ws.Range("E46").Value = "'  -3.62%  "
$ws.Range("D47").Value = "'45.75"
$ws.Range("E47").Value = "'  -1.95%  "
$ws.Range("D48").Value = "'1.32"
$ws.Range("E48").Value = "'  +0.84%  "
$ws.Range("D49").Value = "'7.49"
$ws.Range("E49").Value = "'  -2.22%  "
$ws.Range("D50").Value = "'2.45"
$ws.Range("E50").Value = "'  -3.31%  "
$ws.Range("D51").Value = "'0.987"
$ws.Range("E51").Value = "'  -0.96%  "
